$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 contains a numeric-looking value ("123") in column B that must be
# stored as text (matching the existing data in the sheet, e.g. A2/B2).
# Format the cell as Text before assigning so Excel doesn't coerce it to a number.
$ws.Range("B4").NumberFormat = "@"

$ws.Range("A4").Value = "igna"
$ws.Range("B4").Value = "123"
$ws.Range("C4").Value = "Cliente"
